$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 87: set resultado ("Fallo") and profit (-1)
$ws.Range("G87").Value = "Fallo"
$ws.Range("H87").Value = -1

# Add new row 89 with the latest tracked match
$ws.Range("A89").Value = 14601346
$ws.Range("B89").NumberFormat = "@"
$ws.Range("B89").Value = "2025-09-12"
$ws.Range("C89").Value = "Daniel Evans"
$ws.Range("D89").Value = "Juan Manuel Cerundolo"
$ws.Range("E89").Value = "Gana Juan Manuel Cerundolo"
$ws.Range("F89").Value = 1.83
